$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'27.113.41"
$ws.Range("E2").Value = "`'  +0.68%  "
$ws.Range("D3").Value = "`'1.890.88"
$ws.Range("E3").Value = "`'  +1.57%  "
$ws.Range("D4").Value = "`'1.0000"
$ws.Range("E4").Value = "`'  +0.04%  "
$ws.Range("D5").Value = "`'307.42"
$ws.Range("E5").Value = "`'  +0.93%  "
$ws.Range("D6").Value = "`'0.9993"
$ws.Range("E6").Value = "`'  +0.00%  "
$ws.Range("E7").Value = "`'  +1.80%  "
$ws.Range("D8").Value = "`'0.3731"
$ws.Range("E8").Value = "`'  +2.85%  "
$ws.Range("E9").Value = "`'  +0.93%  "
$ws.Range("D10").Value = "`'21.20"
$ws.Range("E10").Value = "`'  +2.50%  "
$ws.Range("D11").Value = "`'0.9069"
$ws.Range("E11").Value = "`'  +1.48%  "
$ws.Range("D12").Value = "`'0.07630"
$ws.Range("E12").Value = "`'  +2.56%  "
$ws.Range("D13").Value = "`'1.874.97"
$ws.Range("E13").Value = "`'  +1.03%  "
$ws.Range("D14").Value = "`'94.72"
$ws.Range("E14").Value = "`'  +1.81%  "
$ws.Range("E15").Value = "`'  +0.85%  "
$ws.Range("D16").Value = "`'1.001"
$ws.Range("E16").Value = "`'  +0.12%  "
$ws.Range("D17").Value = "`'0.000008518"
$ws.Range("E17").Value = "`'  +0.49%  "
$ws.Range("D18").Value = "`'14.46"
$ws.Range("E18").Value = "`'  +2.59%  "
$ws.Range("D19").Value = "`'0.9993"
$ws.Range("E19").Value = "`'  -0.02%  "
$ws.Range("D20").Value = "`'27.149.42"
$ws.Range("E20").Value = "`'  +0.65%  "
$ws.Range("E21").Value = "`'  +1.23%  "
$ws.Range("D22").Value = "`'2.122.32"
$ws.Range("E22").Value = "`'  +1.57%  "
$ws.Range("D23").Value = "`'10.57"
$ws.Range("E23").Value = "`'  +2.49%  "
$ws.Range("D24").Value = "`'6.427"
$ws.Range("E24").Value = "`'  -0.09%  "
$ws.Range("D25").Value = "`'146.54"
$ws.Range("E25").Value = "`'  -0.92%  "
$ws.Range("D26").Value = "`'1.796"
$ws.Range("E26").Value = "`'  +0.07%  "
$ws.Range("D27").Value = "`'2.209"
$ws.Range("E27").Value = "`'  +7.11%  "
$ws.Range("D28").Value = "`'18.08"
$ws.Range("E28").Value = "`'  +1.31%  "
$ws.Range("D29").Value = "`'114.65"
$ws.Range("E29").Value = "`'  +1.34%  "
$ws.Range("D30").Value = "`'4.971"
$ws.Range("E30").Value = "`'  +6.41%  "
$ws.Range("D31").Value = "`'4.852"
$ws.Range("E31").Value = "`'  +3.99%  "
$ws.Range("D32").Value = "`'0.09208"
$ws.Range("E32").Value = "`'  -0.27%  "
$ws.Range("D33").Value = "`'0.05080"
$ws.Range("E33").Value = "`'  -0.07%  "
$ws.Range("D34").Value = "`'1.239"
$ws.Range("E34").Value = "`'  +7.84%  "
$ws.Range("D35").Value = "`'0.7701"
$ws.Range("E35").Value = "`'  +3.24%  "
$ws.Range("D36").Value = "`'2.989"
$ws.Range("E36").Value = "`'  -0.14%  "
$ws.Range("D37").Value = "`'3.296"
$ws.Range("E37").Value = "`'  +0.43%  "
$ws.Range("D38").Value = "`'2.633"
$ws.Range("E38").Value = "`'  +4.05%  "
$ws.Range("D39").Value = "`'0.5620"
$ws.Range("E39").Value = "`'  +5.04%  "
$ws.Range("D40").Value = "`'0.02000"
$ws.Range("E40").Value = "`'  +0.00%  "
$ws.Range("E41").Value = "`'  -0.73%  "
$ws.Range("D42").Value = "`'6.646"
$ws.Range("E42").Value = "`'  +2.28%  "
$ws.Range("D43").Value = "`'8.965"
$ws.Range("E43").Value = "`'  +5.10%  "
$ws.Range("D44").Value = "`'118.07"
$ws.Range("E44").Value = "`'  +0.08%  "
$ws.Range("D45").Value = "`'0.1512"
$ws.Range("E45").Value = "`'  +3.05%  "
$ws.Range("D46").Value = "`'0.4819"
$ws.Range("E46").Value = "`'  +3.62%  "
$ws.Range("D47").Value = "`'10.23"
$ws.Range("E47").Value = "`'  +1.46%  "
$ws.Range("D48").Value = "`'0.9987"
$ws.Range("E48").Value = "`'  -0.02%  "
$ws.Range("E49").Value = "`'  +2.32%  "
$ws.Range("D50").Value = "`'37.42"
$ws.Range("E50").Value = "`'  +1.54%  "
$ws.Range("D51").Value = "`'64.34"
$ws.Range("E51").Value = "`'  +1.98%  "
